$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emparejamiento")

# --- New cell values on row 2 (tutor 1 side) ---
$ws.Range("C2").Value = "gmail"
$ws.Range("D2").Value = 1234
$ws.Range("N2").Value = 123
$ws.Range("O2").Value = "abc"
$ws.Range("P2").Value = "mama"

# --- New cell values on row 2 (tutor 2 side) ---
$ws.Range("AD2").Value = "mañana"
$ws.Range("AH2").Value = 123
$ws.Range("AI2").Value = "abc"
$ws.Range("AJ2").Value = "papa"
$ws.Range("AX2").Value = "mañana"

# --- Fix subject text on row 3 ---
$ws.Range("G3").Value = "Ingles"

# --- Column widths ---
$ws.Columns.Item(29).ColumnWidth = 28.90625
$ws.Columns.Item(30).ColumnWidth = 24.1796875
$ws.Columns.Item(49).ColumnWidth = 47.08984375

# --- View / selection state ---
$ws.Range("AJ2").Select()
$excel.ActiveWindow.ScrollColumn = 41
